# Generate Report for Handback
# Update the "generated at" / handoff / handback timestamps that get refreshed
# each time the handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 15:15:55"

# zh-cn sheet: "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# for the first file row.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-29 15:15:50"
$wsZhCn.Range("K2").Value = "2016-08-29 15:16:19"

# de-de sheet: "Correspond Handoff Datetime" (mirrors Overview's date for this
# file) and "Correspond Handback DateTime" for the first file row.
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-29 15:15:55"
$wsDeDe.Range("K2").Value = "2016-08-29 15:16:26"
